$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume(1h) updates (text-valued columns; source feed recorded them as plain text)
$ws.Range('D2').Value = '64.641.64'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '3.438.74'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.49'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.99'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('D8').Value = '3.440.98'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  +8.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.34'
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.439'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '4.034.69'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('E15').Value = '  +4.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.32'
$ws.Range('E16').Value = '  +3.44%  '
$ws.Range('D17').Value = '64.669.83'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '3.416.07'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.37'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.25'
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '385.47'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.16'
$ws.Range('E22').Value = '  -3.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.23'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.545'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  +14.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.83'
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.180'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +6.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.44'
$ws.Range('E31').Value = '  +3.82%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.59'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.11'
$ws.Range('E36').Value = '  +3.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.23'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.50'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '3.006.86'
$ws.Range('E39').Value = '  +4.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0765'
$ws.Range('E41').Value = '  -2.63%  '
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.55'
$ws.Range('E43').Value = '  +3.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0317'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.81'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.771'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.72'
$ws.Range('E47').Value = '  +9.37%  '
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('E49').Value = '  +6.79%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.19'
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.62'
$ws.Range('E51').Value = '  +3.81%  '
